$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.050210818428238
$ws.Range("D2").Value = 1.060202080594585
$ws.Range("E2").Value = 1.057551122551956
$ws.Range("F2").Value = 1.068852511993655
$ws.Range("I2").Value = 1.039436020315019
$ws.Range("J2").Value = 1.055245174887831
$ws.Range("K2").Value = 1.062929502774378
$ws.Range("L2").Value = 1.060285782905895
$ws.Range("M2").Value = 1.071556585193251
$ws.Range("N2").Value = 1.056743744074485
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.051299081216518
$ws.Range("D3").Value = 1.061234156569738
$ws.Range("E3").Value = 1.058533324322686
$ws.Range("F3").Value = 1.069955003521541
$ws.Range("I3").Value = 1.039619867726666
$ws.Range("J3").Value = 1.055982350358465
$ws.Range("K3").Value = 1.063775408021388
$ws.Range("L3").Value = 1.061081415605678
$ws.Range("M3").Value = 1.072474423361711
$ws.Range("N3").Value = 1.057481966418842
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.05200300486937
$ws.Range("D4").Value = 1.061902051603675
$ws.Range("E4").Value = 1.059169050044918
$ws.Range("F4").Value = 1.070668715849095
$ws.Range("I4").Value = 1.039736596109165
$ws.Range("J4").Value = 1.05645851704668
$ws.Range("K4").Value = 1.064322230432205
$ws.Range("L4").Value = 1.06159580086049
$ws.Range("M4").Value = 1.073068056723037
$ws.Range("N4").Value = 1.057958809318371
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.052298873619003
$ws.Range("D5").Value = 1.062182851942739
$ws.Range("E5").Value = 1.059436350899425
$ws.Range("F5").Value = 1.070968838574396
$ws.Range("I5").Value = 1.03978513368665
$ws.Range("J5").Value = 1.056658497429004
$ws.Range("K5").Value = 1.064551986369305
$ws.Range("L5").Value = 1.061811942357468
$ws.Range("M5").Value = 1.073317555565682
$ws.Range("N5").Value = 1.058159073695777
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.052348547730699
$ws.Range("D6").Value = 1.062230000562473
$ws.Range("E6").Value = 1.059481234330429
$ws.Range("F6").Value = 1.071019235074505
$ws.Range("I6").Value = 1.039793251972835
$ws.Range("J6").Value = 1.056692063276747
$ws.Range("K6").Value = 1.064590555879317
$ws.Range("L6").Value = 1.061848227240975
$ws.Range("M6").Value = 1.073359443725215
$ws.Range("N6").Value = 1.058192687210874
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.052006958520412
$ws.Range("D7").Value = 1.061905803603945
$ws.Range("E7").Value = 1.059172621568264
$ws.Range("F7").Value = 1.070672725797078
$ws.Range("I7").Value = 1.0397372467735
$ws.Range("J7").Value = 1.056461189980791
$ws.Range("K7").Value = 1.064325300946261
$ws.Range("L7").Value = 1.061598689370751
$ws.Range("M7").Value = 1.073071390792915
$ws.Range("N7").Value = 1.057961486048355
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050578654961054
$ws.Range("D8").Value = 1.060550860311912
$ws.Range("E8").Value = 1.057883025751161
$ws.Range("F8").Value = 1.069225036474263
$ws.Range("I8").Value = 1.039498614564091
$ws.Range("J8").Value = 1.055494479616423
$ws.Range("K8").Value = 1.063215491063445
$ws.Range("L8").Value = 1.060554762181242
$ws.Range("M8").Value = 1.071866828574561
$ws.Range("N8").Value = 1.056993402844387
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.048059831901262
$ws.Range("D9").Value = 1.058163841966133
$ws.Range("E9").Value = 1.055611947530156
$ws.Range("F9").Value = 1.066676528587947
$ws.Range("I9").Value = 1.039061027765743
$ws.Range("J9").Value = 1.053784622437176
$ws.Range("K9").Value = 1.061255774066336
$ws.Range("L9").Value = 1.058711848424793
$ws.Range("M9").Value = 1.069742179089864
$ws.Range("N9").Value = 1.055281117471822
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.046379266536783
$ws.Range("D10").Value = 1.056572881541527
$ws.Range("E10").Value = 1.054098816381255
$ws.Range("F10").Value = 1.064979210005636
$ws.Range("I10").Value = 1.038757836695275
$ws.Range("J10").Value = 1.05264042638499
$ws.Range("K10").Value = 1.05994655082565
$ws.Range("L10").Value = 1.057480972336878
$ws.Range("M10").Value = 1.068324367131198
$ws.Range("N10").Value = 1.054135296530001
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.045651236081648
$ws.Range("D11").Value = 1.055884069386065
$ws.Range("E11").Value = 1.053443833431306
$ws.Range("F11").Value = 1.064244651550685
$ws.Range("I11").Value = 1.038623835825981
$ws.Range("J11").Value = 1.052143959046029
$ws.Range("K11").Value = 1.059378991469379
$ws.Range("L11").Value = 1.056947451977429
$ws.Range("M11").Value = 1.067710111532131
$ws.Range("N11").Value = 1.053638124150473
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.045380762163585
$ws.Range("D12").Value = 1.05562822639268
$ws.Range("E12").Value = 1.053200575542069
$ws.Range("F12").Value = 1.063971862477444
$ws.Range("I12").Value = 1.038573653950448
$ws.Range("J12").Value = 1.051959395169657
$ws.Range("K12").Value = 1.059168075801737
$ws.Range("L12").Value = 1.056749197100555
$ws.Range("M12").Value = 1.067481899369457
$ws.Range("N12").Value = 1.053453298172227
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.045438782056388
$ws.Range("D13").Value = 1.055683105019471
$ws.Range("E13").Value = 1.053252753740268
$ws.Range("F13").Value = 1.064030374009304
$ws.Range("I13").Value = 1.038584436591699
$ws.Range("J13").Value = 1.051998991709561
$ws.Range("K13").Value = 1.059213322403984
$ws.Range("L13").Value = 1.056791727141589
$ws.Range("M13").Value = 1.067530853931034
$ws.Range("N13").Value = 1.053492950943759
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.04562887966765
$ws.Range("D14").Value = 1.055862921070056
$ws.Range("E14").Value = 1.05342372500692
$ws.Range("F14").Value = 1.064222101525702
$ws.Range("I14").Value = 1.03861969609798
$ws.Range("J14").Value = 1.052128706076677
$ws.Range("K14").Value = 1.059361559131876
$ws.Range("L14").Value = 1.056931065833899
$ws.Range("M14").Value = 1.067691248470262
$ws.Range("N14").Value = 1.053622849520155
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.045745998301209
$ws.Range("D15").Value = 1.055973713326898
$ws.Range("E15").Value = 1.053529070278496
$ws.Range("F15").Value = 1.064340238940893
$ws.Range("I15").Value = 1.038641366588165
$ws.Range("J15").Value = 1.052208606995359
$ws.Range("K15").Value = 1.059452879569228
$ws.Range("L15").Value = 1.057016906177145
$ws.Range("M15").Value = 1.067790066171763
$ws.Range("N15").Value = 1.053702863907307
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.046427576409644
$ws.Range("D16").Value = 1.056618597498097
$ws.Range("E16").Value = 1.054142289923851
$ws.Range("F16").Value = 1.065027968466358
$ws.Range("I16").Value = 1.038766672664798
$ws.Range("J16").Value = 1.052673353725664
$ws.Range("K16").Value = 1.059984204016047
$ws.Range("L16").Value = 1.057516368879841
$ws.Range("M16").Value = 1.068365126209425
$ws.Range("N16").Value = 1.054168270631276
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.046855022267725
$ws.Range("D17").Value = 1.057023138749745
$ws.Range("E17").Value = 1.054527003365023
$ws.Range("F17").Value = 1.065459467754164
$ws.Range("I17").Value = 1.03884454656916
$ws.Range("J17").Value = 1.052964603128968
$ws.Range("K17").Value = 1.060317313798563
$ws.Range("L17").Value = 1.057829523224666
$ws.Range("M17").Value = 1.068725756740976
$ws.Range("N17").Value = 1.054459933642139
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.047104311973793
$ws.Range("D18").Value = 1.057259108952306
$ws.Range("E18").Value = 1.054751420743454
$ws.Range("F18").Value = 1.065711191795156
$ws.Range("I18").Value = 1.038889706920864
$ws.Range("J18").Value = 1.053134385292325
$ws.Range("K18").Value = 1.060511547712552
$ws.Range("L18").Value = 1.058012128481773
$ws.Range("M18").Value = 1.068936074245692
$ws.Range("N18").Value = 1.054629956915643
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.047189307876373
$ws.Range("D19").Value = 1.057339570077288
$ws.Range("E19").Value = 1.054827944754323
$ws.Range("F19").Value = 1.065797029606759
$ws.Range("I19").Value = 1.038905060956392
$ws.Range("J19").Value = 1.053192259887189
$ws.Range("K19").Value = 1.06057776571037
$ws.Range("L19").Value = 1.058074383274775
$ws.Range("M19").Value = 1.069007781625468
$ws.Range("N19").Value = 1.054687913699069
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.046809164711316
$ws.Range("D20").Value = 1.056979734461405
$ws.Range("E20").Value = 1.054485725129512
$ws.Range("F20").Value = 1.065413168055008
$ws.Range("I20").Value = 1.038836218552869
$ws.Range("J20").Value = 1.052933365032641
$ws.Range("K20").Value = 1.060281580828596
$ws.Range("L20").Value = 1.057795930144977
$ws.Range("M20").Value = 1.068687067827693
$ws.Range("N20").Value = 1.054428651184132
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.045572902072687
$ws.Range("D21").Value = 1.055809969380782
$ws.Range("E21").Value = 1.053373377348804
$ws.Range("F21").Value = 1.064165640933652
$ws.Range("I21").Value = 1.038609324312797
$ws.Range("J21").Value = 1.052090512675312
$ws.Range("K21").Value = 1.059317909840167
$ws.Range("L21").Value = 1.056890036316988
$ws.Range("M21").Value = 1.067644017665539
$ws.Range("N21").Value = 1.053584601879779
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.044795319298037
$ws.Range("D22").Value = 1.055074563952352
$ws.Range("E22").Value = 1.052674184931183
$ws.Range("F22").Value = 1.063381610302031
$ws.Range("I22").Value = 1.038464306524302
$ws.Range("J22").Value = 1.051559688030377
$ws.Range("K22").Value = 1.058711439639191
$ws.Range("L22").Value = 1.056319992396624
$ws.Range("M22").Value = 1.066987919622891
$ws.Range("N22").Value = 1.053053023402962
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04520755888352
$ws.Range("D23").Value = 1.055464409480088
$ws.Range("E23").Value = 1.053044822559139
$ws.Range("F23").Value = 1.063797207694063
$ws.Range("I23").Value = 1.038541406858838
$ws.Range("J23").Value = 1.051841172615302
$ws.Range("K23").Value = 1.059032995235095
$ws.Range("L23").Value = 1.056622228288166
$ws.Range("M23").Value = 1.067335757247066
$ws.Range("N23").Value = 1.053334907728285
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.046829885866417
$ws.Range("D24").Value = 1.056999346965659
$ws.Range("E24").Value = 1.054504376923717
$ws.Range("F24").Value = 1.065434088779083
$ws.Range("I24").Value = 1.038839982434989
$ws.Range("J24").Value = 1.052947480488493
$ws.Range("K24").Value = 1.060297727215285
$ws.Range("L24").Value = 1.057811109573352
$ws.Range("M24").Value = 1.06870454978151
$ws.Range("N24").Value = 1.054442786685551
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.048711244067414
$ws.Range("D25").Value = 1.058780875894403
$ws.Range("E25").Value = 1.056198914079027
$ws.Range("F25").Value = 1.067335081785196
$ws.Range("I25").Value = 1.039176177023482
$ws.Range("J25").Value = 1.054227418195777
$ws.Range("K25").Value = 1.061762892061825
$ws.Range("L25").Value = 1.059188685772316
$ws.Range("M25").Value = 1.070291695217539
$ws.Range("N25").Value = 1.05572454205119
